$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-14 from 45188 to 45189,
# reflecting the automatic update of the underlying source data by one day.
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45189
}
